$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 3): pitch types
$ws.Range("B3").Value = "Fastball"
$ws.Range("C3").Value = "Riseball"
$ws.Range("D3").Value = "Dropball"
$ws.Range("E3").Value = "Curveball"
$ws.Range("F3").Value = "ChangeUp"
$ws.Range("G3").Value = "Screwball"

# Row labels (column A, rows 4-9): characteristics
$ws.Range("A4").Value = "Velocidad"
$ws.Range("A5").Value = "Posición Cuerpo"
$ws.Range("A6").Value = "Traza"
$ws.Range("A7").Value = "Dirección"
$ws.Range("A8").Value = "Efecto"
$ws.Range("A9").Value = "Agarre"

# Update selection to A10
$ws.Range("A10").Select()
